$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) data range stays text, since many new values
# are numeric-looking strings (e.g. "1.000") that Excel would otherwise
# auto-convert to numbers, losing the original text formatting.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.597.90"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3
$ws.Range("D3").Value = "1.960.02"
$ws.Range("E3").Value = "  +2.15%  "

# Row 4
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "247.87"
$ws.Range("E5").Value = "  +1.01%  "

# Row 6
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("D7").Value = "0.4852"
$ws.Range("E7").Value = "  +1.34%  "

# Row 8
$ws.Range("D8").Value = "44.53"
$ws.Range("E8").Value = "  +1.40%  "

# Row 9
$ws.Range("D9").Value = "0.2934"
$ws.Range("E9").Value = "  +1.12%  "

# Row 10
$ws.Range("D10").Value = "0.06750"
$ws.Range("E10").Value = "  +0.42%  "

# Row 11
$ws.Range("D11").Value = "19.38"
$ws.Range("E11").Value = "  +1.95%  "

# Row 12
$ws.Range("D12").Value = "108.73"
$ws.Range("E12").Value = "  -1.79%  "

# Row 13
$ws.Range("D13").Value = "1.967.62"
$ws.Range("E13").Value = "  +2.74%  "

# Row 14
$ws.Range("D14").Value = "0.07762"
$ws.Range("E14").Value = "  +2.49%  "

# Row 15
$ws.Range("D15").Value = "5.428"
$ws.Range("E15").Value = "  +2.65%  "

# Row 16
$ws.Range("D16").Value = "0.6829"
$ws.Range("E16").Value = "  +1.86%  "

# Row 17
$ws.Range("D17").Value = "291.81"
$ws.Range("E17").Value = "  -2.33%  "

# Row 18
$ws.Range("D18").Value = "30.625.22"
$ws.Range("E18").Value = "  +0.26%  "

# Row 19
$ws.Range("D19").Value = "13.18"
$ws.Range("E19").Value = "  +1.73%  "

# Row 20
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.223.60"
$ws.Range("E20").Value = "  +2.67%  "

# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.000007659"
$ws.Range("E21").Value = "  +1.08%  "

# Row 22
$ws.Range("D22").Value = "5.602"
$ws.Range("E22").Value = "  -0.42%  "

# Row 23
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.01%  "

# Row 25
$ws.Range("D25").Value = "6.612"
$ws.Range("E25").Value = "  +1.43%  "

# Row 26
$ws.Range("D26").Value = "9.845"
$ws.Range("E26").Value = "  +3.86%  "

# Row 27
$ws.Range("D27").Value = "169.56"
$ws.Range("E27").Value = "  +2.96%  "

# Row 28
$ws.Range("D28").Value = "19.94"
$ws.Range("E28").Value = "  -1.48%  "

# Row 29
$ws.Range("D29").Value = "2.163"
$ws.Range("E29").Value = "  +2.63%  "

# Row 30
$ws.Range("D30").Value = "0.1062"
$ws.Range("E30").Value = "  -0.90%  "

# Row 31
$ws.Range("D31").Value = "1.438"
$ws.Range("E31").Value = "  +2.79%  "

# Row 32
$ws.Range("D32").Value = "4.653"
$ws.Range("E32").Value = "  +14.97%  "

# Row 33
$ws.Range("D33").Value = "4.418"
$ws.Range("E33").Value = "  +6.35%  "

# Row 34
$ws.Range("D34").Value = "0.05084"
$ws.Range("E34").Value = "  +1.47%  "

# Row 35
$ws.Range("D35").Value = "0.7641"
$ws.Range("E35").Value = "  +3.35%  "

# Row 36
$ws.Range("D36").Value = "1.172"
$ws.Range("E36").Value = "  +2.84%  "

# Row 37
$ws.Range("D37").Value = "2.733"
$ws.Range("E37").Value = "  -0.03%  "

# Row 38
$ws.Range("D38").Value = "0.02027"
$ws.Range("E38").Value = "  +0.03%  "

# Row 39
$ws.Range("D39").Value = "2.720"
$ws.Range("E39").Value = "  +1.32%  "

# Row 40
$ws.Range("D40").Value = "6.487"
$ws.Range("E40").Value = "  +10.48%  "

# Row 41
$ws.Range("D41").Value = "2.117"
$ws.Range("E41").Value = "  +4.45%  "

# Row 42
$ws.Range("D42").Value = "0.4442"
$ws.Range("E42").Value = "  -0.71%  "

# Row 43
$ws.Range("D43").Value = "108.80"
$ws.Range("E43").Value = "  -2.09%  "

# Row 44
$ws.Range("D44").Value = "0.8734"
$ws.Range("E44").Value = "  +1.02%  "

# Row 45
$ws.Range("D45").Value = "70.06"
$ws.Range("E45").Value = "  -2.09%  "

# Row 46
$ws.Range("E46").Value = "  +0.15%  "

# Row 47
$ws.Range("D47").Value = "7.466"
$ws.Range("E47").Value = "  +2.80%  "

# Row 48
$ws.Range("D48").Value = "0.1273"
$ws.Range("E48").Value = "  +2.87%  "

# Row 49
$ws.Range("D49").Value = "9.348"
$ws.Range("E49").Value = "  +0.19%  "

# Row 50
$ws.Range("D50").Value = "35.81"
$ws.Range("E50").Value = "  +2.10%  "

# Row 51
$ws.Range("D51").Value = "47.29"
$ws.Range("E51").Value = "  -4.48%  "
